$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 45; this shifts the existing rows
# 45-79 down to 46-80 (the former last row, row 79, becomes row 80).
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new weekly price record.
$ws.Cells.Item(45, 1).Value = 6
$ws.Cells.Item(45, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(45, 3).Value = "Metropolitana"
$ws.Cells.Item(45, 4).Value = 45090
$ws.Cells.Item(45, 5).Value = 13
$ws.Cells.Item(45, 6).Value = 100112035
$ws.Cells.Item(45, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 500
$ws.Cells.Item(45, 11).Value = 18000
$ws.Cells.Item(45, 12).Value = 20000
$ws.Cells.Item(45, 13).Value = 19080
$ws.Cells.Item(45, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(45, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(45, 16).Value = 1272
$ws.Cells.Item(45, 17).Value = 15
$ws.Cells.Item(45, 18).Value = "Hortaliza"
